$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 101
$ws.Range("H101").Value = 101410.6
$ws.Range("I101").Value = 1786.7142
$ws.Range("J101").Value = 333866.34
$ws.Range("K101").Value = 5360.142599999999
$ws.Range("L101").Value = 1001599.02
$ws.Range("M101").Value = -3738.142599999999
$ws.Range("N101").Value = -1004843.02
# Row 111
$ws.Range("H111").Value = 1104.2858
$ws.Range("I111").Value = 1121.6666
$ws.Range("K111").Value = 3364.9998
$ws.Range("M111").Value = -297.9998000000001
# Row 132
$ws.Range("H132").Value = 2085087.1
$ws.Range("I132").Value = 1567.025
$ws.Range("K132").Value = 4701.075000000001
$ws.Range("M132").Value = -2171.075000000001
# Row 137
$ws.Range("H137").Value = 814.8511
$ws.Range("I137").Value = 754.94116
$ws.Range("K137").Value = 2264.82348
$ws.Range("M137").Value = 285.17652
# Row 138
$ws.Range("H138").Value = 2526.27
$ws.Range("I138").Value = 1128.025
$ws.Range("J138").Value = 3458.4333
$ws.Range("K138").Value = 3384.075
$ws.Range("L138").Value = 10375.2999
$ws.Range("M138").Value = 1755.925
$ws.Range("N138").Value = -20655.2999
# Row 139
$ws.Range("H139").Value = 91944.445
$ws.Range("J139").Value = 91944.445
$ws.Range("L139").Value = 91944.445
$ws.Range("N139").Value = -102224.445

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1496.6154
$ws.Range("I45").Value = 1511.7778
$ws.Range("J45").Value = 1462.5
$ws.Range("K45").Value = 1511.7778
$ws.Range("L45").Value = 1462.5
$ws.Range("M45").Value = -1134.7778
$ws.Range("N45").Value = -2216.5
# Row 74
$ws.Range("H74").Value = 1157.4054
$ws.Range("I74").Value = 1215.4166
$ws.Range("J74").Value = 1050.3077
$ws.Range("K74").Value = 1215.4166
$ws.Range("L74").Value = 1050.3077
$ws.Range("M74").Value = -341.4166
$ws.Range("N74").Value = -2798.3077
# Row 77
$ws.Range("H77").Value = 1157.4054
$ws.Range("I77").Value = 1215.4166
$ws.Range("J77").Value = 1050.3077
$ws.Range("K77").Value = 6077.083000000001
$ws.Range("L77").Value = 5251.538500000001
$ws.Range("M77").Value = -1709.083000000001
$ws.Range("N77").Value = -13987.5385
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()
# Row 122
$ws.Range("H122").Value = 752.125
$ws.Range("I122").Value = 759.5714
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 2278.7142
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = 171.2857999999997
$ws.Range("N122").Value = -7000

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 47400
$ws.Range("J20").Value = 47400
$ws.Range("L20").Value = 47400
$ws.Range("N20").Value = -47872
# Row 30
$ws.Range("H30").Value = 47400
$ws.Range("J30").Value = 47400
$ws.Range("L30").Value = 47400
$ws.Range("N30").Value = -47582
# Row 31
$ws.Range("H31").Value = 24271.873
$ws.Range("I31").Value = 2917
$ws.Range("J31").Value = 65646.94
$ws.Range("K31").Value = 2917
$ws.Range("L31").Value = 65646.94
$ws.Range("M31").Value = -2622
$ws.Range("N31").Value = -66236.94
# Row 34
$ws.Range("H34").Value = 24271.873
$ws.Range("I34").Value = 2917
$ws.Range("J34").Value = 65646.94
$ws.Range("K34").Value = 2917
$ws.Range("L34").Value = 65646.94
$ws.Range("M34").Value = -2715
$ws.Range("N34").Value = -66050.94
# Row 50
$ws.Range("H50").Value = 9215
$ws.Range("J50").Value = 9215
$ws.Range("L50").Value = 9215
$ws.Range("N50").Value = -10465
# Row 51
$ws.Range("H51").Value = 8607
$ws.Range("J51").Value = 9079.2
$ws.Range("L51").Value = 9079.2
$ws.Range("N51").Value = -10551.2
# Row 59
$ws.Range("H59").Value = 16027
$ws.Range("J59").Value = 16027
$ws.Range("L59").Value = 16027
$ws.Range("N59").Value = -18317
# Row 60
$ws.Range("H60").Value = 8145.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 8145.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 8145.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -9167.5
# Row 61
$ws.Range("H61").Value = 8607
$ws.Range("J61").Value = 9079.2
$ws.Range("L61").Value = 9079.2
$ws.Range("N61").Value = -9775.2
# Row 68
$ws.Range("H68").Value = 14462.375
$ws.Range("J68").Value = 14462.375
$ws.Range("L68").Value = 14462.375
$ws.Range("N68").Value = -15960.375
# Row 71
$ws.Range("H71").Value = 14462.375
$ws.Range("J71").Value = 14462.375
$ws.Range("L71").Value = 43387.125
$ws.Range("N71").Value = -50875.125
# Row 74
$ws.Range("H74").Value = 13714
$ws.Range("J74").Value = 13714
$ws.Range("L74").Value = 13714
$ws.Range("N74").Value = -15462
# Row 77
$ws.Range("H77").Value = 13714
$ws.Range("J77").Value = 13714
$ws.Range("L77").Value = 41142
$ws.Range("N77").Value = -49878
# Row 99
$ws.Range("H99").Value = 2176.5
$ws.Range("I99").Value = 2087.6428
$ws.Range("J99").Value = 2487.5
$ws.Range("K99").Value = 2087.6428
$ws.Range("L99").Value = 2487.5
$ws.Range("M99").Value = -589.6428000000001
$ws.Range("N99").Value = -5483.5
# Row 126
$ws.Range("H126").Value = 2176.5
$ws.Range("I126").Value = 2087.6428
$ws.Range("J126").Value = 2487.5
$ws.Range("K126").Value = 6262.928400000001
$ws.Range("L126").Value = 7462.5
$ws.Range("M126").Value = -3792.928400000001
$ws.Range("N126").Value = -12402.5
# Row 128
$ws.Range("H128").Value = 47400
$ws.Range("J128").Value = 47400
$ws.Range("L128").Value = 47400
$ws.Range("N128").Value = -57360
# Row 134
$ws.Range("H134").Value = 20001192
$ws.Range("I134").Value = 1264.4546
$ws.Range("K134").Value = 3793.3638
$ws.Range("M134").Value = -1258.3638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 781.6667
$ws.Range("I113").Value = 1305.25
$ws.Range("J113").Value = 591.2727
$ws.Range("K113").Value = 3915.75
$ws.Range("L113").Value = 1773.8181
$ws.Range("M113").Value = -1745.75
$ws.Range("N113").Value = -6113.8181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 41
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -645
# Row 104
$ws.Range("H104").Value = 35000
$ws.Range("J104").Value = 35000
$ws.Range("L104").Value = 35000
$ws.Range("N104").Value = -41988
# Row 122
$ws.Range("H122").Value = 1013185.75
$ws.Range("I122").Value = 1196855.9
$ws.Range("K122").Value = 3590567.7
$ws.Range("M122").Value = -3588117.7
# Row 132
$ws.Range("H132").Value = 3714.8064
$ws.Range("I132").Value = 3929.4443
$ws.Range("J132").Value = 3417.6155
$ws.Range("K132").Value = 11788.3329
$ws.Range("L132").Value = 10252.8465
$ws.Range("M132").Value = -9258.332900000001
$ws.Range("N132").Value = -15312.8465

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2619.25
$ws.Range("I7").Value = 2579.7778
$ws.Range("J7").Value = 2974.5
$ws.Range("K7").Value = 2579.7778
$ws.Range("L7").Value = 2974.5
$ws.Range("M7").Value = -2467.7778
$ws.Range("N7").Value = -3198.5
# Row 122
$ws.Range("H122").Value = 2508.8
$ws.Range("I122").Value = 2508.8
$ws.Range("K122").Value = 7526.400000000001
$ws.Range("M122").Value = -5076.400000000001
# Row 126
$ws.Range("H126").Value = 2619.25
$ws.Range("I126").Value = 2579.7778
$ws.Range("J126").Value = 2974.5
$ws.Range("K126").Value = 7739.3334
$ws.Range("L126").Value = 8923.5
$ws.Range("M126").Value = -5269.3334
$ws.Range("N126").Value = -13863.5
# Row 136
$ws.Range("H136").Value = 2694.92
$ws.Range("I136").Value = 1692.28
$ws.Range("J136").Value = 4700.2
$ws.Range("K136").Value = 5076.84
$ws.Range("L136").Value = 14100.6
$ws.Range("M136").Value = -2526.84
$ws.Range("N136").Value = -19200.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 18531.25
$ws.Range("J86").Value = 18531.25
$ws.Range("L86").Value = 18531.25
$ws.Range("N86").Value = -20777.25
# Row 89
$ws.Range("H89").Value = 18531.25
$ws.Range("J89").Value = 18531.25
$ws.Range("L89").Value = 92656.25
$ws.Range("N89").Value = -103888.25
